$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $searchText) {
    $r = $doc.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return -1
    }
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $r.Start -and $p.Range.End -ge $r.End) {
            return $i
        }
    }
    return -1
}

# --- 1. "Cite (MLA):" paragraph gains a hanging indent ---
$idx1 = Find-ParagraphIndex $d "Cite (MLA):"
if ($idx1 -eq -1) { throw "Could not find 'Cite (MLA):' paragraph" }
$p1 = $d.Paragraphs.Item($idx1)
$xml1 = @"
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" mc:Ignorable="w14" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006"><w:body><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:left="720" w:hanging="720"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Cite (MLA):</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$p1.Range.InsertXML($xml1)

# --- 2. Citation paragraph: reformatted/merged runs (Times-Roman/Times-Italic,
#        direct hyperlink formatting, ligatures) + new trailing empty paragraph ---
$idx2 = Find-ParagraphIndex $d "Steere, Edward, ed."
if ($idx2 -eq -1) { throw "Could not find citation paragraph" }
$p2 = $d.Paragraphs.Item($idx2)
$xml2 = @"
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" mc:Ignorable="w14" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006"><w:body><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:left="720" w:hanging="720"/><w:rPr><w:rFonts w:ascii="Times-Roman" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Times-Roman" w:cs="Times-Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times-Roman" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Times-Roman" w:cs="Times-Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr><w:t xml:space="preserve">Steere, Edward, ed. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times-Italic" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Times-Italic" w:cs="Times-Italic"/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr><w:t>Swahili Tales, as Told by Natives of Zanzibar, with an English Translation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times-Roman" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Times-Roman" w:cs="Times-Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr><w:t xml:space="preserve">. Translated by Edward Steere. London: Bell &amp; Daldy, 1870. </w:t></w:r><w:hyperlink r:id="rId9" w:history="1"><w:r><w:rPr><w:rFonts w:ascii="Times-Roman" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Times-Roman" w:cs="Times-Roman"/><w:color w:val="0000E9"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single" w:color="0000E9"/><w:lang w:val="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr><w:t>https://archive.org/details/swahilitalesasto00stee</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rFonts w:ascii="Times-Roman" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Times-Roman" w:cs="Times-Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$p2.Range.InsertXML($xml2)

# --- 3. Drop the stale <w:lastRenderedPageBreak/> cache marker on the
#        "The Story of the Washerman's Donkey" Heading 1 run. A no-op
#        find/replace over just that run rewrites it (and so clears the
#        render cache marker) without disturbing the surrounding bookmarks. ---
$idx3 = Find-ParagraphIndex $d "“The Story of the Washerman’s Donkey”"
if ($idx3 -eq -1) { throw "Could not find the Washerman's Donkey heading paragraph" }
$p3 = $d.Paragraphs.Item($idx3)
$headingText = "“The Story of the Washerman’s Donkey”"
$p3.Range.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, $headingText, 2) | Out-Null

Write-Output "done"
